$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$oldValues = @(
    "48+38="
    "53-49="
    "23+49="
    "8+46="
    "84-49="
    "70-34="
    "28+46="
    "50-1="
    "17+24="
    "36-19="
    "19+44="
    "84-36="
    "92-25="
    "37+46="
    "54-7="
    "33-28="
    "5+77="
    "16+38="
    "16+37="
    "6+19="
    "36-9="
    "84-75="
    "8+79="
    "9+34="
    "71-66="
    "40-21="
    "8+29="
    "50-5="
    "44+28="
    "2+39="
    "60-51="
    "8+39="
    "30-29="
    "41-18="
    "38+46="
    "57+37="
    "27+35="
    "51-36="
    "76-8="
    "5+67="
    "45+17="
    "92-53="
    "60-43="
    "55+36="
    "49+15="
    "12-6="
    "40-1="
    "74-26="
    "19+67="
    "67+28="
    "92-15="
    "31-2="
    "37+38="
    "96-17="
    "89+5="
    "41-15="
    "26+28="
    "28+55="
    "38+5="
    "61-56="
    "19+55="
    "26+36="
    "7+68="
    "33+8="
    "2+19="
    "66+9="
    "14+78="
    "29+57="
    "35-27="
    "18+3="
    "6+37="
    "45-38="
    "36+49="
    "44+18="
    "54+28="
    "7+46="
    "22-17="
    "86-78="
    "67+15="
    "89+8="
    "60-34="
    "72-38="
    "4+59="
    "9+37="
    "17-9="
    "96-78="
    "19+57="
    "35+57="
    "64-57="
    "6+76="
    "18+6="
    "29+26="
    "27+65="
    "62-24="
    "59+33="
    "20-11="
    "9+58="
    "8+45="
    "83-25="
    "92-74="
)

$newValues = @(
    "29+52="
    "17+57="
    "76-38="
    "18+74="
    "77+8="
    "39+8="
    "7+88="
    "48+37="
    "66-49="
    "12+69="
    "3+89="
    "42+39="
    "40-36="
    "9+69="
    "83-18="
    "47+26="
    "61-29="
    "60-4="
    "26+49="
    "72-26="
    "45-19="
    "18+68="
    "22+19="
    "65-29="
    "43-18="
    "64-56="
    "71-46="
    "20-1="
    "54-36="
    "60-19="
    "65-18="
    "59+22="
    "24+49="
    "52+19="
    "95-8="
    "22-3="
    "48+14="
    "18+7="
    "62-23="
    "58+9="
    "63-59="
    "54-35="
    "48-29="
    "9+79="
    "51-24="
    "51-23="
    "17-8="
    "66-48="
    "6+49="
    "28+36="
    "76-59="
    "53-47="
    "26+37="
    "37+54="
    "65-28="
    "34+48="
    "94-48="
    "23+68="
    "9+77="
    "81-54="
    "91-4="
    "47+46="
    "66+8="
    "91-26="
    "22+29="
    "39+38="
    "43-8="
    "90-26="
    "67-58="
    "24+8="
    "63-8="
    "66-18="
    "51-13="
    "39+12="
    "15+9="
    "72-18="
    "62+9="
    "9+5="
    "9+4="
    "65+29="
    "92-68="
    "38+59="
    "49+42="
    "73-17="
    "18+64="
    "81-34="
    "57+29="
    "66-58="
    "20-9="
    "56-18="
    "54+38="
    "28+16="
    "54-6="
    "2+29="
    "5+57="
    "56+28="
    "62-45="
    "64-48="
    "91-44="
    "22-15="
)

$idx = 0
$mismatches = 0
$rowCount = $tbl.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $row = $tbl.Rows.Item($r)
    $colCount = $row.Cells.Count
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $row.Cells.Item($c)
        $current = $cell.Range.Text
        $current = $current.TrimEnd([char]13, [char]7)
        $expectedOld = $oldValues[$idx]
        if ($current -ne $expectedOld) {
            Write-Output "MISMATCH at idx $idx : expected [$expectedOld] found [$current]"
            $mismatches = $mismatches + 1
        }
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
Write-Output "Updated $idx cells, $mismatches mismatches"
